# Add header row: Datum / Plaats / Duiker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datum"
$ws.Range("B1").Value = "Plaats"
$ws.Range("C1").Value = "Duiker"

# A1 already carries the bold/bordered/centered header style (s="1").
# Copy that formatting onto the two new header cells so all three match.
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
